$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139, shifting existing rows 139-210 down to 140-211.
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with its values.
# Columns A,B,C,E,F,G,H,I,J are constant across all data rows in this sheet.
$ws.Cells.Item(139, 1).Value = 11
$ws.Cells.Item(139, 2).Value = "Vega Monumental Concepci$([char]0xF3)n"
$ws.Cells.Item(139, 3).Value = "B$([char]0xED)ob$([char]0xED)o"
$ws.Cells.Item(139, 4).Value = 44845
$ws.Cells.Item(139, 5).Value = 8
$ws.Cells.Item(139, 6).Value = "Fruta"
$ws.Cells.Item(139, 7).Value = 100108
$ws.Cells.Item(139, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(139, 9).Value = 100108005
$ws.Cells.Item(139, 10).Value = "Pi$([char]0xF1)a"
$ws.Cells.Item(139, 11).Value = "Caramelo"
$ws.Cells.Item(139, 12).Value = "Segunda"
$ws.Cells.Item(139, 13).Value = 200
$ws.Cells.Item(139, 14).Value = 20000
$ws.Cells.Item(139, 15).Value = 21000
$ws.Cells.Item(139, 16).Value = 20500
$ws.Cells.Item(139, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(139, 18).Value = "Ecuador"
$ws.Cells.Item(139, 19).Value = 1464
$ws.Cells.Item(139, 20).Value = 14

# Ensure the date cell keeps the date-formatted style used by the rest of column D.
$ws.Cells.Item(139, 4).NumberFormat = $ws.Cells.Item(140, 4).NumberFormat
